$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new value in E5 (Drivetrain Instantiation (multiple threads) for the 1,000,000 cohort row)
$ws.Range("E5").Value = "20 hours 33 minutes"

# Update H5 with new measurement value
$ws.Range("H5").Value = 3522298786

# Update the view: scroll back to A1 (no frozen/top-left offset) and move selection to G7
$ws.Range("A1").Select()
$ws.Range("G7").Select()
